$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.118.62'
$ws.Range("E2").Value = '  -0.91%  '

$ws.Range("D3").Value = '2.240.97'
$ws.Range("E3").Value = '  -1.57%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.630'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.58'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +6.15%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.627'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.32'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +6.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0948'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.16'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.87%  '

$ws.Range("E13").Value = '  -2.09%  '

$ws.Range("D14").Value = '2.578.12'
$ws.Range("E14").Value = '  -1.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.87'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.862'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.30%  '

$ws.Range("D17").Value = '2.235.80'
$ws.Range("E17").Value = '  -2.30%  '

$ws.Range("D18").Value = '42.039.21'
$ws.Range("E18").Value = '  -1.12%  '

$ws.Range("D19").Value = '0.0₃0979'
$ws.Range("E19").Value = '  -1.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.15'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.45'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.52'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.83%  '

$ws.Range("E23").Value = '  -3.00%  '

$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.72'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.19'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.31'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.31'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +14.17%  '

$ws.Range("E29").Value = '  -1.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.40'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.50'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0853'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.67'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.25%  '

$ws.Range("E34").Value = '  -5.59%  '

$ws.Range("E35").Value = '  +1.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.61'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.85'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0297'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.17'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.20'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.88'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '114.57'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +18.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.203'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.98%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.15'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.82'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.100'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.85%  '

$ws.Range("E47").Value = '  -0.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.13'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.11%  '

$ws.Range("E49").Value = '  -1.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.29'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -12.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.27'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.08%  '
